# ------------------------------------------------------------------
# Refresh the crypto price/volume snapshot (cryptos.xlsx) to match the
# latest scrape. The sheet stores every data cell as plain text (the
# "Price" column even holds values such as "34.704.81" that LOOK like
# numbers but are not meant to be parsed as one). To stop Excel from
# "helpfully" auto-converting number-like strings into real numbers -
# which would also silently attach a numeric style to the cell - each
# cell is explicitly forced to the Text format before its value is
# written, and the style is reset back to Normal immediately after so
# no stray formatting is left behind.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "34.704.81"),
    @("E2", "  +0.67%  "),
    @("D3", "1.819.39"),
    @("E3", "  +1.17%  "),
    @("E4", "  +0.22%  "),
    @("D5", "228.51"),
    @("E5", "  +0.61%  "),
    @("D6", "0.577"),
    @("E6", "  +3.82%  "),
    @("E7", "  +0.11%  "),
    @("D8", "34.75"),
    @("E8", "  +6.95%  "),
    @("E9", "  +1.19%  "),
    @("D10", "0.0702"),
    @("E10", "  +1.01%  "),
    @("E11", "  +0.44%  "),
    @("D12", "2.083.69"),
    @("E12", "  +1.35%  "),
    @("D13", "11.42"),
    @("E13", "  +2.89%  "),
    @("D14", "1.820.99"),
    @("E14", "  +1.01%  "),
    @("E15", "  +1.53%  "),
    @("D16", "34.681.29"),
    @("E16", "  +0.73%  "),
    @("D17", "4.34"),
    @("E17", "  +2.07%  "),
    @("D18", "69.16"),
    @("E18", "  +0.92%  "),
    @("D19", "0.0₃0803"),
    @("E19", "  +0.02%  "),
    @("D20", "246.85"),
    @("E20", "  -0.04%  "),
    @("D21", "11.61"),
    @("E21", "  +4.15%  "),
    @("E22", "  +0.14%  "),
    @("D23", "4.19"),
    @("E23", "  +0.54%  "),
    @("D24", "173.83"),
    @("E24", "  +6.46%  "),
    @("E25", "  +1.27%  "),
    @("D26", "7.50"),
    @("E26", "  +3.12%  "),
    @("D27", "16.83"),
    @("E27", "  +1.98%  "),
    @("E28", "  +2.41%  "),
    @("E29", "  -0.01%  "),
    @("E30", "  +2.44%  "),
    @("E31", "  +1.66%  "),
    @("E32", "  +1.94%  "),
    @("E33", "  +0.94%  "),
    @("E34", "  +0.95%  "),
    @("D35", "2.64"),
    @("E35", "  +1.18%  "),
    @("D36", "1.410.99"),
    @("E36", "  -2.39%  "),
    @("E37", "  +1.90%  "),
    @("E38", "  +1.66%  "),
    @("E39", "  +0.46%  "),
    @("D40", "84.82"),
    @("E40", "  +0.76%  "),
    @("D41", "2.86"),
    @("E41", "  +4.16%  "),
    @("D42", "0.954"),
    @("E42", "  +2.06%  "),
    @("E43", "  -0.05%  "),
    @("E44", "  -0.27%  "),
    @("E45", "  +3.09%  "),
    @("D46", "0.0518"),
    @("E46", "  -1.33%  "),
    @("E47", "  -0.09%  "),
    @("D48", "1.983.95"),
    @("E48", "  +1.68%  "),
    @("D49", "105.61"),
    @("E49", "  -0.23%  "),
    @("B50", "BabyDogeCoin"),
    @("C50", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"),
    @("D50", "0.0₆0131"),
    @("E50", "  +1.11%  "),
    @("B51", "PaxDollar"),
    @("C51", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"),
    @("D51", "1.00"),
    @("E51", "  +0.15%  ")
)

foreach ($pair in $updates) {
    $addr = $pair[0]
    $newValue = $pair[1]
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}
